$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 447
$col = 3  # Column C ("Förändrad")

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = 46075
}
